$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 18: Pre-Alpha Dev. entry
$ws.Range("A18").Value = 45956
$ws.Range("B18").Value = "Pre-Alpha Dev."
$ws.Range("C18").Value = 0.45833333333333331
$ws.Range("D18").Value = 0.5
$ws.Range("F18").Value = "Record and edit microprocessor portion of video"

# Row 19: Pre-Alpha Dev. entry
$ws.Range("A19").Value = 45956
$ws.Range("B19").Value = "Pre-Alpha Dev."
$ws.Range("C19").Value = 0.65972222222222221
$ws.Range("D19").Value = 0.68402777777777779
$ws.Range("F19").Value = "Clean up readme and edit in frontend portion of video"

# Update selection to F20 as per diff
$ws.Range("F20").Select()
